# Practical 8 - update the "count" column (D) on the autophagosome sheet
# with the actual counts, after the tutorial. Cells that already read 0
# stay as-is; the ones below get their real tallies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$counts = @{
    3  = "1"
    4  = "10"
    5  = "3"
    7  = "2"
    15 = "1"
    24 = "1"
    28 = "2"
    33 = "5"
    36 = "2"
}

foreach ($row in $counts.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $counts[$row]
}
